$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.140.36"
$ws.Range("E2").Value = "  +0.94%  "
$ws.Range("D3").Value = "2.050.87"
$ws.Range("E3").Value = "  -3.30%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.52%  "
$ws.Range("E6").Value = "  -2.70%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "55.91"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +16.46%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "61.87"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.377"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0760"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.61%  "
$ws.Range("E12").Value = "  +5.94%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.08"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.71%  "
$ws.Range("D14").Value = "2.347.13"
$ws.Range("E14").Value = "  -3.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.822"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.27%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.04%  "
$ws.Range("D17").Value = "2.050.47"
$ws.Range("E17").Value = "  -3.24%  "
$ws.Range("D18").Value = "37.026.57"
$ws.Range("E18").Value = "  +0.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.09%  "
$ws.Range("D20").Value = "0.0₃0888"
$ws.Range("E20").Value = "  +5.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.41"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.90%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("E25").Value = "  -2.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.22"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.27%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.99"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.15%  "
$ws.Range("E30").Value = "  -0.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.59"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.96%  "
$ws.Range("E32").Value = "  +14.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0623"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.78%  "
$ws.Range("E34").Value = "  +4.30%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0865"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -9.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.27"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.21%  "
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.78"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.08%  "
$ws.Range("B39").Value = "Gas"
$ws.Range("C39").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.56"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -31.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.111"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +31.79%  "
$ws.Range("E41").Value = "  -0.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "18.26"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +12.53%  "
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.14"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.16%  "
$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +64.36%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "97.03"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.20%  "
$ws.Range("D48").Value = "1.301.89"
$ws.Range("E48").Value = "  -4.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.89%  "
$ws.Range("E50").Value = "  +3.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.84"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.35%  "
